# Fruta / hortaliza, semanal
# Insert the latest weekly price entries (2 rows) right after the existing
# most-recent week (rows 2-3), shifting all historical rows down by two
# and appending the corresponding data at the bottom via the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4; existing rows 4:24 move to 6:26.
$ws.Rows("4:5").Insert()

# New row 4: Primera
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 45133
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 100112043
$ws.Range("G4").Value = "Pepino dulce"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 68
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("N4").Value = "`$/caja 15 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 1200
$ws.Range("Q4").Value = 15
$ws.Range("R4").Value = "Hortaliza"

# New row 5: Segunda
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 45133
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 100112043
$ws.Range("G5").Value = "Pepino dulce"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 13000
$ws.Range("N5").Value = "`$/caja 15 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 867
$ws.Range("Q5").Value = 15
$ws.Range("R5").Value = "Hortaliza"

Write-Host "Inserted two new rows with latest weekly data; dimension now" $ws.UsedRange.Address()
